$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C7) from 45208 to 45212
$ws.Range("C2:C7").Value = 45212
